# Applies the "Master's in Marketing" / "Bachelor's in Business" education
# block restructuring + skills tweak + new "Fluent in Spanish" line.

$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$apos = [char]0x2019   # curly RIGHT SINGLE QUOTATION MARK used in "Master's"

# ---------------------------------------------------------------------------
# 1) "Masters Studies in Marketing" -> "Master's" / " in " / "Marketing"
#    (three bold runs) + the _GoBack bookmark moves here.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(90)
$xml1 = "<w:p $ns>" +
  "<w:pPr><w:ind w:left=`"100`" w:right=`"6426`"/><w:jc w:val=`"both`"/>" +
    "<w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/></w:rPr>" +
  "</w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/></w:rPr>" +
    "<w:t>Master${apos}s</w:t></w:r>" +
  "<w:r><w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/></w:rPr>" +
    "<w:t xml:space=`"preserve`"> in </w:t></w:r>" +
  "<w:r><w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/></w:rPr>" +
    "<w:t>Marketing</w:t></w:r>" +
  "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" +
  "</w:p>"
$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) First "{UNIVERSITY}" line keeps its text but the right-indent narrows
#    from 5312 to 4722 (same as the Bachelor's one further down).
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(91)
$xml2 = "<w:p $ns>" +
  "<w:pPr><w:spacing w:before=`"39`"/><w:ind w:left=`"100`" w:right=`"4722`"/><w:jc w:val=`"both`"/>" +
    "<w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/></w:rPr>" +
  "</w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/></w:rPr>" +
    "<w:t>{UNIVERSITY}</w:t></w:r>" +
  "</w:p>"
$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) The old blank spacer paragraph becomes an empty (no run) paragraph
#    carrying the bold/sz21 "Bachelor's" formatting.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(92)
$xml3 = "<w:p $ns>" +
  "<w:pPr><w:ind w:left=`"100`" w:right=`"6426`"/><w:jc w:val=`"both`"/>" +
    "<w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/></w:rPr>" +
  "</w:pPr>" +
  "</w:p>"
$p3.Range.InsertXML($xml3)

# ---------------------------------------------------------------------------
# 4) "Bachelor's in Business" keeps its text, moves to the ind right=6426
#    paragraph style (matching the one vacated above).
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(93)
$xml4 = "<w:p $ns>" +
  "<w:pPr><w:ind w:left=`"100`" w:right=`"6426`"/><w:jc w:val=`"both`"/>" +
    "<w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/></w:rPr>" +
  "</w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/></w:rPr>" +
    "<w:t>Bachelor's in Business</w:t></w:r>" +
  "</w:p>"
$p4.Range.InsertXML($xml4)

# ---------------------------------------------------------------------------
# 5) Second "{UNIVERSITY}" line keeps its text/indent but loses the
#    _GoBack bookmark (it now lives on the Master's paragraph).
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(94)
$xml5 = "<w:p $ns>" +
  "<w:pPr><w:spacing w:before=`"39`"/><w:ind w:left=`"100`" w:right=`"4722`"/><w:jc w:val=`"both`"/>" +
    "<w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/></w:rPr>" +
  "</w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/></w:rPr>" +
    "<w:t>{UNIVERSITY}</w:t></w:r>" +
  "</w:p>"
$p5.Range.InsertXML($xml5)

# ---------------------------------------------------------------------------
# 6) The spacer paragraph right before "SKILLS" changes its pPr shape
#    (ind right=8123 / jc both / sz18, same look as the SKILLS heading).
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(98)
$xml6 = "<w:p $ns>" +
  "<w:pPr><w:ind w:left=`"100`" w:right=`"8123`"/><w:jc w:val=`"both`"/>" +
    "<w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/></w:rPr>" +
  "</w:pPr>" +
  "</w:p>"
$p6.Range.InsertXML($xml6)

# ---------------------------------------------------------------------------
# 7) Drop the "Social Media, " prefix from the skills list.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Social Media, Marketing, Sales Coordination", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Marketing, Sales Coordination", 2) | Out-Null

# ---------------------------------------------------------------------------
# 8) Append a blank line followed by "Fluent in Spanish", matching the
#    formatting of the skills paragraph above them.
# ---------------------------------------------------------------------------
$xmlBlank = "<w:p $ns>" +
  "<w:pPr><w:spacing w:line=`"301`" w:lineRule=`"auto`"/><w:ind w:left=`"100`" w:right=`"89`"/><w:jc w:val=`"both`"/>" +
    "<w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/></w:rPr>" +
  "</w:pPr>" +
  "</w:p>"
$xmlSpanish = "<w:p $ns>" +
  "<w:pPr><w:spacing w:line=`"301`" w:lineRule=`"auto`"/><w:ind w:left=`"100`" w:right=`"89`"/><w:jc w:val=`"both`"/>" +
    "<w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/></w:rPr>" +
  "</w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:eastAsia=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/></w:rPr>" +
    "<w:t>Fluent in Spanish</w:t></w:r>" +
  "</w:p>"

$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertXML($xmlBlank)

$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertXML($xmlSpanish)
